$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add Priority column label ---
$ws.Cells.Item(1,1).Value = 'Priority'

# --- Data rows 2-21 (re-sorted/updated backlog, by priority) ---
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 'Display message to registered users who have not yet paid'
$ws.Cells.Item(2,3).Value = 4

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 'Auto-pick team at deadline'
$ws.Cells.Item(3,3).Value = 6

$ws.Cells.Item(4,1).Value = 1
$ws.Cells.Item(4,2).Value = 'Competition Standings'
$ws.Cells.Item(4,3).Value = 4

$ws.Cells.Item(5,1).Value = 1
$ws.Cells.Item(5,2).Value = 'Include information on charity'
$ws.Cells.Item(5,3).Value = 2

$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = 'Splash screen for first time users summerising competition rules'
$ws.Cells.Item(6,3).Value = 5

$ws.Cells.Item(7,1).Value = 1
$ws.Cells.Item(7,2).Value = 'Payment page'
$ws.Cells.Item(7,3).Value = 5

$ws.Cells.Item(8,1).Value = 1
$ws.Cells.Item(8,2).Value = 'Refactor headers and footers'
$ws.Cells.Item(8,3).Value = 1
$ws.Cells.Item(8,4).Value = 1

$ws.Cells.Item(9,1).Value = 1
$ws.Cells.Item(9,2).Value = 'Refactor DAL - include disconnects, monitor performance'
$ws.Cells.Item(9,3).Value = 3

$ws.Cells.Item(10,1).Value = 1
$ws.Cells.Item(10,2).Value = 'Add "Full name" field to registration a'
$ws.Cells.Item(10,3).Value = 1

$ws.Cells.Item(11,1).Value = 2
$ws.Cells.Item(11,2).Value = 'Where a user is eliminated, indicate this at login, show history'
$ws.Cells.Item(11,3).Value = 4

$ws.Cells.Item(12,1).Value = 2
$ws.Cells.Item(12,2).Value = 'Player History'
$ws.Cells.Item(12,3).Value = 4

$ws.Cells.Item(13,1).Value = 2
$ws.Cells.Item(13,2).Value = 'Header tabs to select "Home|User Settings | league"'
$ws.Cells.Item(13,3).Value = 2

$ws.Cells.Item(14,1).Value = 2
$ws.Cells.Item(14,2).Value = 'Fix footer - change to navbar'
$ws.Cells.Item(14,3).Value = 2

$ws.Cells.Item(15,1).Value = 2
$ws.Cells.Item(15,2).Value = 'Fix "History for... " injection of username'
$ws.Cells.Item(15,3).Value = 1
$ws.Cells.Item(15,5).Value = 'in progress'

$ws.Cells.Item(16,1).Value = 3
$ws.Cells.Item(16,2).Value = 'Disable teams that user cannot select'
$ws.Cells.Item(16,3).Value = 2
$ws.Cells.Item(16,5).Value = 'Done'

$ws.Cells.Item(17,1).Value = 3
$ws.Cells.Item(17,2).Value = 'Image background'
$ws.Cells.Item(17,3).Value = 3

$ws.Cells.Item(18,1).Value = 3
$ws.Cells.Item(18,2).Value = 'Style text content blocks correctly'
$ws.Cells.Item(18,3).Value = 3

$ws.Cells.Item(19,1).Value = 3
$ws.Cells.Item(19,2).Value = 'Style confirmation dialogs (avoid browser native)'
$ws.Cells.Item(19,3).Value = 2

$ws.Cells.Item(20,1).Value = 3
$ws.Cells.Item(20,2).Value = 'Show team/match selection in confirm dialog before submitting prediction'
$ws.Cells.Item(20,3).Value = 1

$ws.Cells.Item(21,1).Value = 3
$ws.Cells.Item(21,2).Value = 'Fix header navbar (white highlight on middle 2 buttons)'
$ws.Cells.Item(21,3).Value = 2

# --- Clear stale D16 (Act Effort moved to row 8) ---
$ws.Cells.Item(16,4).ClearContents()

# --- Trailing blank-ish rows 22-24 (single-space placeholder) ---
$ws.Cells.Item(22,1).Value = ' '
$ws.Cells.Item(23,1).Value = ' '
$ws.Cells.Item(24,1).Value = ' '

# --- Apply "Good" conditional style to completed/near-complete tasks ---
$ws.Range("B2").Style = "Good"
$ws.Range("B4").Style = "Good"
$ws.Range("B5").Style = "Good"
$ws.Range("B7").Style = "Good"
$ws.Range("B8").Style = "Good"
$ws.Range("B9").Style = "Good"
$ws.Range("B12").Style = "Good"
$ws.Range("B13").Style = "Good"

# --- Re-apply sort by Priority (col A) over the data block, records sortState ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A21"))
$sortObj.SetRange($ws.Range("A2:D21"))
$sortObj.Header = 0
$sortObj.Apply()

# --- Update selection to match final cursor position ---
$ws.Range("B24").Select()

Write-Host "edit complete"
